$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.815.01"
$ws.Range("E2").Value = "  -1.74%  "

# Row 3
$ws.Range("D3").Value = "1.889.78"
$ws.Range("E3").Value = "  -1.74%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.11%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7679"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.15%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "244.47"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.01%  "

# Row 7
$ws.Range("E7").Value = "  -0.10%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3121"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.30%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "25.21"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -7.22%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07206"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.91%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08105"
$ws.Range("D11").Style = "Normal"

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7647"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.59%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.498"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.45%  "

# Row 14
$ws.Range("D14").Value = "1.908.60"
$ws.Range("E14").Value = "  -1.01%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.27"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.25%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.136"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.55%  "

# Row 17
$ws.Range("D17").Value = "29.822.46"
$ws.Range("E17").Value = "  -1.83%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.90"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.10%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "242.67"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.45%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007760"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.35%  "

# Row 21
$ws.Range("E21").Value = "  -0.06%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.176"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.90%  "

# Row 23
$ws.Range("D23").Value = "2.150.11"
$ws.Range("E23").Value = "  -2.17%  "

# Row 24
$ws.Range("E24").Value = "  -0.15%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1558"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -6.83%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.397"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.38%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.42"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.14%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.72"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.06%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.039"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.60%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.471"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.92%  "

# Row 31
$ws.Range("E31").Value = "  +0.11%  "

# Row 32
$ws.Range("E32").Value = "  +2.32%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.085"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.40%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05517"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.37%  "

# Row 35
$ws.Range("E35").Value = "  -3.58%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7470"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.04%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.002"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.16%  "

# Row 38
$ws.Range("E38").Value = "  -3.48%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01922"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.01%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.781"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.49%  "

# Row 41
$ws.Range("D41").Value = "1.152.21"
$ws.Range("E41").Value = "  +11.15%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "73.61"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.95%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4413"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.10%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.902"
$ws.Range("D44").Style = "Normal"

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8489"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.96%  "

# Row 46
$ws.Range("E46").Value = "  -0.07%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.82"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.40%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.882"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.48%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.885"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.58%  "

# Row 50
$ws.Range("E50").Value = "  -0.76%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.440"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.79%  "
